# Update "想去人数" (want-to-go count) values in the "展览" and "全部类型"
# sheets. Each listed cell value is incremented by 1.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")
$updates = @{
    "F2" = 339
    "F3" = 89
    "F5" = 20
    "F6" = 46
    "F7" = 128
    "F8" = 50
    "F9" = 334
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
